$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D1").Value = "ORG_FIE_IDENOLD"
$ws.Range("E1").Value = "ORG_FIE_IDENNEW"
$ws.Range("F1").Value = "ORG_FIE_STATUS"

$ws.Range("D10").Select()
